$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 9-10 for the two new sale items. This pushes the
# existing running-total row (old row 9) down to row 11 and the footer
# row (old row 10) down to row 12, preserving their content/format.
$ws.Rows("9:10").Insert()

# The new rows are additional sale-item lines - clone the layout
# (styles + merges) of the existing item row (row 8) onto them.
$ws.Range("A8:Q8").Copy($ws.Range("A9"))
$ws.Range("A8:Q8").Copy($ws.Range("A10"))
$ws.Rows(9).RowHeight = 25.5
$ws.Rows(10).RowHeight = 24.75

# ---- Row 9: new item "حنه جلوري سوده 1 كيس" ----
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "حنه جلوري سوده 1 كيس"
$ws.Range("H9").Value = "3:0"
$ws.Range("L9").Value = "'0"
$ws.Range("N9").Value = "45.00"
$ws.Range("P9").Value = "'45.0000"
$ws.Range("Q9").Value = "1:0"

# ---- Row 10: new item "سرنجه دواء" ----
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "سرنجه دواء"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = "'0"
$ws.Range("N10").Value = "4.00"
$ws.Range("P10").Value = "'4.0000"
$ws.Range("Q10").Value = "1:0"

# ---- Row 11 (previously row 9): refreshed running total ----
$ws.Range("P11").Value2 = 197

# ---- Row 12 (previously row 10): footer, refreshed generation timestamp ----
$ws.Range("A12").Value = "Sunday, 27 July, 2025 9:58 AM"

Write-Output "edit applied"
